$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 60, pushing existing row 60 (and below) down to row 61
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new SBE-911 instrument entry
$ws.Range("A60").Value = "SBE-911"
$ws.Range("B60").Value = "S1460"

# Update the active selection / view to reflect where the edit was made
$excel.ActiveWindow.ScrollRow = 46
$ws.Range("B60").Select()
